$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "mutual_fund_instrument_name"
$ws.Range("D1").Select()
